$d = $word.ActiveDocument

# 1. "Begin " -> "BEGIN" (drop trailing space, upper-case) in the top-level
#    intro paragraph. Find/Replace keeps a clean <w:t> (no xml:space) and
#    only matches the body text, not the "Begin" caption inside the shape.
$d.Content.Find.Execute("Begin ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "BEGIN", 2)

# 2. "End" -> "END" in the closing paragraph (again, only the body text).
$d.Content.Find.Execute("End", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "END", 2)

# Locate the paragraph that now just contains "END".
$endParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "END`r") {
        $endParagraph = $para
        break
    }
}

$endPos = $endParagraph.Range.End - 1

# 3. Move the "_GoBack" bookmark so it sits right after the "END" run
#    (collapsed, zero-length) instead of its old spot before the second
#    flowchart shape. A temporary character is inserted/removed around the
#    target position so the bookmark lands exactly at the run boundary;
#    re-adding a bookmark with the existing name relocates it (bookmark
#    names are unique), so the old occurrence disappears on its own.
$tempRange = $d.Range($endPos, $endPos)
$tempRange.InsertAfter("X")

$markRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markRange)

$charRange = $d.Range($endPos, $endPos + 1)
$charRange.Delete()

Write-Host "done"
